$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.139.14'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '3.083.60'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('D5').Value = '''577.15'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').Value = '''169.25'
$ws.Range('E6').Value = '  -2.21%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '3.080.82'
$ws.Range('E8').Value = '  -0.84%  '
$ws.Range('E9').Value = '  -1.04%  '
$ws.Range('D10').Value = '''6.40'
$ws.Range('E10').Value = '  -0.97%  '
$ws.Range('D11').Value = '''0.151'
$ws.Range('E11').Value = '  -1.00%  '
$ws.Range('E12').Value = '  -1.26%  '
$ws.Range('D13').Value = '''0.0000242'
$ws.Range('E13').Value = '  -1.42%  '
$ws.Range('D14').Value = '''36.18'
$ws.Range('E14').Value = '  -2.44%  '
$ws.Range('E15').Value = '  -2.06%  '
$ws.Range('D16').Value = '3.596.78'
$ws.Range('E16').Value = '  -0.80%  '
$ws.Range('D17').Value = '67.026.97'
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').Value = '''7.03'
$ws.Range('E18').Value = '  -0.96%  '
$ws.Range('D19').Value = '''16.59'
$ws.Range('E19').Value = '  +1.82%  '
$ws.Range('D20').Value = '3.085.10'
$ws.Range('E20').Value = '  -0.85%  '
$ws.Range('D21').Value = '''491.73'
$ws.Range('E21').Value = '  +3.21%  '
$ws.Range('B22').Value = 'Polygon'
$ws.Range('C22').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D22').Value = '''0.689'
$ws.Range('E22').Value = '  -3.15%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = '''7.70'
$ws.Range('E23').Value = '  -1.53%  '
$ws.Range('D24').Value = '''82.87'
$ws.Range('E24').Value = '  -1.04%  '
$ws.Range('D25').Value = '''12.92'
$ws.Range('E25').Value = '  -3.51%  '
$ws.Range('D26').Value = '''2.22'
$ws.Range('E26').Value = '  -2.46%  '
$ws.Range('D27').Value = '''10.24'
$ws.Range('E27').Value = '  +3.39%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').Value = '''7.88'
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('D30').Value = '''2.30'
$ws.Range('E30').Value = '  -5.11%  '
$ws.Range('D31').Value = '''2.62'
$ws.Range('E31').Value = '  -1.03%  '
$ws.Range('D32').Value = '''27.90'
$ws.Range('E32').Value = '  -2.61%  '
$ws.Range('E33').Value = '  -1.59%  '
$ws.Range('D34').Value = '0.0₃0910'
$ws.Range('E34').Value = '  -3.35%  '
$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').Value = '''5.71'
$ws.Range('E36').Value = '  -2.50%  '
$ws.Range('D38').Value = '''46.92'
$ws.Range('E38').Value = '  -1.17%  '
$ws.Range('E39').Value = '  +1.10%  '
$ws.Range('D40').Value = '''2.00'
$ws.Range('E40').Value = '  -3.91%  '
$ws.Range('D41').Value = '''0.304'
$ws.Range('E41').Value = '  -1.93%  '
$ws.Range('D42').Value = '''8.32'
$ws.Range('E42').Value = '  -3.09%  '
$ws.Range('D43').Value = '2.775.86'
$ws.Range('E43').Value = '  -0.48%  '
$ws.Range('D44').Value = '''371.95'
$ws.Range('E44').Value = '  -1.48%  '
$ws.Range('E45').Value = '  -2.61%  '
$ws.Range('D46').Value = '''135.86'
$ws.Range('D47').Value = '''2.46'
$ws.Range('E47').Value = '  -3.29%  '
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('E49').Value = '  -0.81%  '
$ws.Range('E50').Value = '  -1.70%  '
$ws.Range('E51').Value = '  -1.31%  '
